$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.919.04"
$ws.Range("E2").Value = "  +4.94%  "
$ws.Range("D3").Value = "4.037.43"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.24"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.87"
$ws.Range("E6").Value = "  +9.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.696"
$ws.Range("E7").Value = "  +14.76%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.760"
$ws.Range("E9").Value = "  +7.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000330"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.68"
$ws.Range("E12").Value = "  +17.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.86"
$ws.Range("E13").Value = "  +5.61%  "
$ws.Range("D14").Value = "4.680.07"
$ws.Range("E14").Value = "  +4.95%  "
$ws.Range("D15").Value = "4.054.41"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.33"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.73"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "71.885.53"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.19"
$ws.Range("E21").Value = "  +5.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "99.75"
$ws.Range("E22").Value = "  +15.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.70"
$ws.Range("E24").Value = "  +5.25%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.24"
$ws.Range("E25").Value = "  +5.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.19"
$ws.Range("E26").Value = "  -6.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.71"
$ws.Range("E28").Value = "  +30.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.13"
$ws.Range("E29").Value = "  +5.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.83"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.62"
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("E32").Value = "  +6.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "680.67"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.87"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.77"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.73"
$ws.Range("E36").Value = "  +8.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.430"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.157"
$ws.Range("E38").Value = "  +7.00%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").Value = "  +13.28%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0844"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0492"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("E45").Value = "  +7.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  -6.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.59"
$ws.Range("E48").Value = "  +10.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.06"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.37"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("E51").Value = "  -2.91%  "
